# BAU Subsidy per Unit Electricity Supplied to the Grid by Grid Batteries
# Update the "About" sheet Notes section: replace the US Inflation Reduction
# Act / battery-ITC note with notes describing the (lack of) EU support
# scheme for grid batteries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("A6").Value  = "In the EU there was no general support for grid batteries"
$ws.Range("A7").Value  = "in place looking at the time before the Green Deal"
$ws.Range("A8").Value  = "There were some projects financed under PCI (Projets of Common Interest)"
$ws.Range("A9").Value  = "and IPCEI (Important Project of Common European Interest)"
$ws.Range("A10").Value = "but there was no financial support applicable to any grid battery"

$wb.Save()
